# Daily attendance processing - 2026-01-04 17:31:45
# Normalizes the "Recorded By" (column G) values: when the leading entry
# in the comma-separated list is "System" (any casing), move it to the
# end of the list and normalize its casing to "System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 157
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Length -gt 1 -and $parts[0].ToLower() -eq "system") {
        $lastIdx = $parts.Length - 1
        $parts[0] = $parts[$lastIdx]
        $parts[$lastIdx] = "System"
        $newVal = $parts -join ", "
        $cell.Value = $newVal
    }
}
